$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 12.6251
$ws.Range("A8").Value = -21.10660000000001
$ws.Range("A10").Value = -20.49339999999997
$ws.Range("A12").Value = -22.31740000000003
$ws.Range("D13").Value = -7.661500000000001
$ws.Range("A18").Value = -22.25090000000002
$ws.Range("E20").Value = 12.06709999999999
$ws.Range("A25").Value = -22.19830000000003
